$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "CaixaDeTexto 16" textbox on slide 1 (shape id 17) that holds
# the DISCIPLINA / QUALIDADE DE SOFTWARE ... text block.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "CaixaDeTexto 16") {
        $shp = $candidate
    }
}

$tr = $shp.TextFrame.TextRange

# The shape currently has 3 paragraphs:
#   1) "DISCIPLINA:   PROJETO DE SISTEMAS APLICADO AS MELHORES PRATICAS EM "
#   2) "QUALIDADE DE SOFTWARE E GOVERNANCA DE TI"
#   3) "" (empty, trailing)
# Insert a new paragraph "TURMA 3SI" right after paragraph 2, pushing the
# existing empty paragraph down to become paragraph 4.
$p2 = $tr.Paragraphs(2)
[void]$p2.InsertAfter("`rTURMA 3SI")
